$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)      # "Sheet1" (original)

# Duplicate Sheet1 -> Excel auto-names the copy "Sheet1 (2)" and places it
# right after the original. This copy keeps the full R1:R5 rule data and
# later gets the data bugfix too.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)      # "Sheet1 (2)" (the duplicate)

# Bugfix (rule data): the Age rule for R4 should be "-" (N/A), not "T".
# A leading apostrophe makes Excel store it as quote-prefixed text,
# matching the style already used elsewhere for the "-" placeholder.
$ws2.Range("H3").Value = "'-"

# The duplicate also picks up a new "X" mark at F7.
$ws2.Range("F7").Value = "X"

# Its lingering selection moves to H4 (it is no longer the active tab).
[void]$ws2.Range("H4").Select()

# Back on the original "Sheet1": trim it down to only show ID / Variable /
# Operator / Value plus the R4 and R5 rule columns -- drop the R1:R3
# columns (old E:G) so old H:I slide left into E:F.
$ws1.Range("E:G").EntireColumn.Delete()

# Apply the same data bugfix in its new column position (was H3, now E3).
$ws1.Range("E3").Value = "'-"

# Selection on "Sheet1" lands on F4.
[void]$ws1.Range("F4").Select()

# "Sheet1" (the trimmed summary view) is the tab left active/selected.
[void]$ws1.Activate()
